$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph that currently sits right after
#    the document title (Heading1) paragraph.
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# 2. Insert a new bold paragraph ("Play Book of Oz Lock 'n Spin Free | Slot
#    Game Review") right before the trailing "Prompt for DALLE" paragraph.
$n = $d.Paragraphs.Count
$dallePara = $d.Paragraphs($n)
$dallePara.Range.InsertParagraphBefore() | Out-Null

$n = $d.Paragraphs.Count
$newPara = $d.Paragraphs($n - 1)
$newParaXml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Book of Oz Lock ''n Spin Free | Slot Game Review</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara.Range.InsertXML($newParaXml)

# 3. Replace the DALLE prompt text with the new meta-description text, while
#    keeping the paragraph's italic run formatting and its leading empty run.
$n = $d.Paragraphs.Count
$dallePara = $d.Paragraphs($n)
$textRange = $d.Range($dallePara.Range.Start, $dallePara.Range.End - 1)
$textRange.Text = "Explore the Wizard of Oz world with Book of Oz Lock 'n Spin. Get the Lock 'n Spin feature, customizable paylines and play for free."
